$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.648.82'
$ws.Range('E2').Value = '  +0.34%  '

$ws.Range('D3').Value = '1.701.99'
$ws.Range('E3').Value = '  +0.94%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.33%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.59%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3967'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.86%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4055'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.85%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.001'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.514'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.13%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08795'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.79%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.337'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +11.42%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.29'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.35%  '

$ws.Range('E15').Value = '  +1.59%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.535'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.62%  '

$ws.Range('D17').Value = '1.704.94'
$ws.Range('E17').Value = '  +0.68%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '100.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07127'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.80%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.53'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.39%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.759'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.37%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.34%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.47%  '

$ws.Range('D24').Value = '24.632.72'
$ws.Range('E24').Value = '  +0.33%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.012'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.19%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.319'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.42%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.73%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.24%  '

$ws.Range('B29').Value = 'HuobiToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.140'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.25%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.51%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.412'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +27.70%  '

$ws.Range('D32').Value = '1.894.64'
$ws.Range('E32').Value = '  +0.76%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.086'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.93%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08676'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.84%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.317'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +18.27%  '

$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.965'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.09%  '

$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.70%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2721'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.43%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.75'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.75%  '

$ws.Range('E40').Value = '  +9.35%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.08996'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.67%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.476'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.22%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7689'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.44%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7187'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.22%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.58%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.458'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.180'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.54%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9995'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.33%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.27%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.296'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.19%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00000000376'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.73%  '
